$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header cells C1 and D1 had their text swapped: C1 ("notes ") <-> D1 ("notes")
# and the "notes " label picked up a ".1" suffix (duplicate-header disambiguation),
# becoming "notes.1" instead of "notes".
$ws.Range("C1").Value = "notes"
$ws.Range("D1").Value = "notes.1"

# Selection moved from P1:P1048576 (whole column P) to the single cell D1.
$ws.Range("D1").Select() | Out-Null
